# The workbook's player list (A2:D164) gets sorted by Price (column D),
# descending — largest price first. This mirrors selecting column D and
# choosing "Sort Largest to Smallest" from the Data ribbon / context menu.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D164")
$keyRange  = $ws.Range("D93:D164")

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($keyRange, 0, 2, $null, 0)
$sortObj.SetRange($dataRange)
$sortObj.Header = 2
$sortObj.MatchCase = $false
$sortObj.Orientation = 1
$sortObj.SortMethod = 1
$sortObj.Apply()

# Reflect the user's final selection: the whole Price column, scrolled so
# row 93 is at the top of the view.
$ws.Columns("D").Select() | Out-Null
